# Daten aktualisiert am 2024-03-19
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @("TAO-USD", "IMX-USD", "GRT-USD", "PEPE-USD", "MNT-USD")

$startRow = 404
for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}
